# Weekly update: insert the latest Cilantro price-sheet entry as a new
# row 28, pushing the existing rows 28-84 down to 29-85.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 28 (shifts 28..84 down to 29..85,
# carrying their formatting/styles with them).
$ws.Rows.Item(28).Insert()

# Populate the newly inserted row 28 with the new week's data.
$ws.Range("A28").Value = 8
$ws.Range("B28").Value = "Terminal La Palmera de La Serena"
$ws.Range("C28").Value = "Coquimbo"
$ws.Range("D28").Value = 44469
$ws.Range("E28").Value = 4
$ws.Range("F28").Value = 100112040
$ws.Range("G28").Value = "Cilantro"
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 3160
$ws.Range("K28").Value = 1500
$ws.Range("L28").Value = 2000
$ws.Range("M28").Value = 1750
$ws.Range("N28").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O28").Value = "Provincia del Elquí"
$ws.Range("P28").Value = 1167
$ws.Range("Q28").Value = 1.5
$ws.Range("R28").Value = "Hortaliza"
